$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename column header from "ID User" to "ID Supplier" (stok: user_id -> supplier_id)
$ws.Range("B1").Value = "ID Supplier"

# Update the active selection to match the saved view state
$ws.Range("F10").Select()
